$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "779÷9="; New = "490÷2=" },
    @{ Row = 1;  Col = 2; Old = "724÷8="; New = "941÷4=" },
    @{ Row = 1;  Col = 3; Old = "245÷6="; New = "434÷8=" },
    @{ Row = 1;  Col = 4; Old = "630÷7="; New = "250÷2=" },
    @{ Row = 1;  Col = 5; Old = "553÷9="; New = "297÷8=" },

    @{ Row = 5;  Col = 1; Old = "859÷2="; New = "760÷8=" },
    @{ Row = 5;  Col = 2; Old = "135÷6="; New = "785÷7=" },
    @{ Row = 5;  Col = 3; Old = "916÷9="; New = "918÷8=" },
    @{ Row = 5;  Col = 4; Old = "403÷7="; New = "542÷8=" },
    @{ Row = 5;  Col = 5; Old = "941÷8="; New = "347÷9=" },

    @{ Row = 9;  Col = 1; Old = "756÷3="; New = "755÷5=" },
    @{ Row = 9;  Col = 2; Old = "316÷7="; New = "163÷6=" },
    @{ Row = 9;  Col = 3; Old = "224÷9="; New = "732÷4=" },
    @{ Row = 9;  Col = 4; Old = "347÷6="; New = "632÷9=" },
    @{ Row = 9;  Col = 5; Old = "827÷5="; New = "275÷7=" },

    @{ Row = 13; Col = 1; Old = "146÷4="; New = "819÷2=" },
    @{ Row = 13; Col = 2; Old = "640÷8="; New = "755÷3=" },
    @{ Row = 13; Col = 3; Old = "734÷3="; New = "688÷6=" },
    @{ Row = 13; Col = 4; Old = "577÷3="; New = "103÷8=" },
    @{ Row = 13; Col = 5; Old = "700÷8="; New = "475÷2=" },

    @{ Row = 17; Col = 1; Old = "640÷8="; New = "668÷4=" },
    @{ Row = 17; Col = 2; Old = "100÷3="; New = "371÷5=" },
    @{ Row = 17; Col = 3; Old = "901÷4="; New = "893÷2=" },
    @{ Row = 17; Col = 4; Old = "689÷2="; New = "707÷4=" },
    @{ Row = 17; Col = 5; Old = "106÷7="; New = "270÷5=" }
)

foreach ($item in $replacements) {
    $cell = $tbl.Cell($item.Row, $item.Col)
    $rng = $cell.Range
    $rng.Find.Execute($item.Old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $item.New, 1)
}
